# Generate Report for Handoff
# Updates the generated UUID / content-hash file names, the handoff
# timestamps, and the matching hyperlink display text across the
# "Overview", "zh-cn" and "de-de" worksheets.

$wb = $excel.ActiveWorkbook

$oldGuid = "63f28594-7c66-4be7-af39-dfd1ae3af9ba"
$newGuid = "7fe0fe36-b965-48b3-9c80-4a0d3ac847a7"
$oldHash = "04434e24def8b289fe4641e3ff7cdb2412bc09db"
$newHash = "e69752e6b1b36c488f0e86371a82f940f49ec798"

$oldMd   = "$oldGuid.md"
$newMd   = "$newGuid.md"
$oldZh   = "$oldGuid.$oldHash.zh-cn.xlf"
$newZh   = "$newGuid.$newHash.zh-cn.xlf"
$oldDe   = "$oldGuid.$oldHash.de-de.xlf"
$newDe   = "$newGuid.$newHash.de-de.xlf"

# External hyperlink targets are unchanged by this edit - only the
# on-sheet display text changes - so the exact same addresses used by
# the workbook today are reused when the hyperlinks are recreated.
$addrMd = "https://github.com/OpenLocalizationTest/oltest/blob/6ca1bcc2b284bf66aff800e26346cde4f5457127/e2e/$oldMd"
$addrZh = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c671264df7b759c13b75c19006790108ff22286b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$oldZh"
$addrDe = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4d4963d3eed9b8e24d231eb7c767a53072ba145f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$oldDe"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# Hyperlinks.Delete() removes every hyperlink on the sheet, so all of
# them are recreated afterwards (here there is only the one on A2).
$ws.Range("A2").Hyperlinks.Delete()

$ws.Range("A2").Value = $newMd
$ws.Range("D2").Value = "2016-47-17 22:47:47"

$ws.Hyperlinks.Add($ws.Range("A2"), $addrMd, "", "", $newMd)

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Hyperlinks.Delete()

$ws.Range("A2").Value = $newMd
$ws.Range("D2").Value = $newZh
$ws.Range("E2").Value = "2016-03-17 22:47:44"

$ws.Hyperlinks.Add($ws.Range("A2"), $addrMd, "", "", $newMd)
$ws.Hyperlinks.Add($ws.Range("B2"), $addrMd, "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D2"), $addrZh, "", "", $newZh)

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Hyperlinks.Delete()

$ws.Range("A2").Value = $newMd
$ws.Range("D2").Value = $newDe
$ws.Range("E2").Value = "2016-03-17 22:47:47"

$ws.Hyperlinks.Add($ws.Range("A2"), $addrMd, "", "", $newMd)
$ws.Hyperlinks.Add($ws.Range("B2"), $addrMd, "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D2"), $addrDe, "", "", $newDe)
